# Generate Report for Handback
# Adds a new handback row (ac481a19-cbc4-42bf-abde-74a8bf5386cb.md) to the
# Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

$fileId = "ac481a19-cbc4-42bf-abde-74a8bf5386cb"

# ---------------------------------------------------------------------------
# Overview sheet (row 4)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "$fileId.md"
$wsOverview.Range("B4").Value = "e2e\$fileId.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G4").Value = "2017-02-09 09:51:06"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c7a89fd2c3843e97e05643dfcb1bb33563e00ba4/e2e/$fileId.md", "", "", "e2e\$fileId.md") | Out-Null

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------------
# zh-cn sheet (row 4)
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A4").Value = "$fileId.md"
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "'True"
$wsZhCn.Range("G4").Value = "$fileId.c7a89fd2c3843e97e05643dfcb1bb33563e00ba4.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2017-02-09 09:50:43"
$wsZhCn.Range("I4").Value = "'"
$wsZhCn.Range("J4").Value = "$fileId.md"
$wsZhCn.Range("K4").Value = "$fileId.c7a89fd2c3843e97e05643dfcb1bb33563e00ba4.zh-cn.xlf"
$wsZhCn.Range("L4").Value = "2017-02-09 09:51:50"
$wsZhCn.Range("M4").Value = "'"
$wsZhCn.Range("N4").Value = "'"
$wsZhCn.Range("O4").Value = "'True"
$wsZhCn.Range("P4").Value = "'"
$wsZhCn.Range("Q4").Value = "'False"
$wsZhCn.Range("R4").Value = "'"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c7a89fd2c3843e97e05643dfcb1bb33563e00ba4/e2e/$fileId.md", "", "", "$fileId.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("J4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c7a89fd2c3843e97e05643dfcb1bb33563e00ba4/e2e/$fileId.md", "", "", "$fileId.md") | Out-Null

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:R4"))

# ---------------------------------------------------------------------------
# de-de sheet (row 4)
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A4").Value = "$fileId.md"
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "'True"
$wsDeDe.Range("G4").Value = "$fileId.c7a89fd2c3843e97e05643dfcb1bb33563e00ba4.de-de.xlf"
$wsDeDe.Range("H4").Value = "2017-02-09 09:51:06"
$wsDeDe.Range("I4").Value = "'"
$wsDeDe.Range("J4").Value = "$fileId.md"
$wsDeDe.Range("K4").Value = "$fileId.c7a89fd2c3843e97e05643dfcb1bb33563e00ba4.de-de.xlf"
$wsDeDe.Range("L4").Value = "2017-02-09 09:52:17"
$wsDeDe.Range("M4").Value = "'"
$wsDeDe.Range("N4").Value = "'"
$wsDeDe.Range("O4").Value = "'True"
$wsDeDe.Range("P4").Value = "'"
$wsDeDe.Range("Q4").Value = "'False"
$wsDeDe.Range("R4").Value = "'"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c7a89fd2c3843e97e05643dfcb1bb33563e00ba4/e2e/$fileId.md", "", "", "$fileId.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("J4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/c7a89fd2c3843e97e05643dfcb1bb33563e00ba4/e2e/$fileId.md", "", "", "$fileId.md") | Out-Null

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:R4"))
